$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# --- Update the three tier input counts (row 11) ---
$ws.Range("D11").Value = 500
$ws.Range("E11").Value = 1000
$ws.Range("I11").Value = 7500
$ws.Range("J11").Value = 10000

# --- Add new informational rows 16 and 17 ---
# Set the new text / formula content first, in the same left-to-right,
# top-to-bottom order the cells were authored so new shared-strings land
# at the expected indices (112 .. 115).
$ws.Range("N17").Value = "Generate Link that they click then we can count the visitors per sended mails"
$ws.Range("C16").Value = "incl Btw"
$ws.Range("C17").Value = "Active"
$ws.Range("D17").Value = "x"
$ws.Range("F17").Value = "x"
$ws.Range("H17").Value = "x"

$ws.Range("D16").Formula = "=D10/79%"
$ws.Range("E16").Formula = "=E10/79%"
$ws.Range("F16").Formula = "=F10/79%"
$ws.Range("H16").Formula = "=H10/79%"
$ws.Range("I16").Formula = "=I10/79%"
$ws.Range("J16").Formula = "=J10/79%"

# G16 stays empty but picks up the same currency number format as the
# rest of row 16 / row 10; F17 keeps its "x" value but also picks up that
# currency format (matches the source workbook). Use a format-only paste
# so we reuse the existing style instead of minting a new one.
$ws.Range("D10").Copy()
$ws.Range("G16").PasteSpecial(-4122)

$ws.Range("D10").Copy()
$ws.Range("F17").PasteSpecial(-4122)

# --- Move the active selection to N27 (matches the saved view state) ---
$ws.Range("N27").Select()

Write-Output "done"
